$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column at M (13) for "Revolution in Earth Days" ---
# This shifts the old M,N,O,P columns one to the right (N,O,P,Q).
$ws.Columns.Item(13).Insert()

# Old "Revolution - one earth year" column (L) is now hidden; the new
# "Revolution in Earth Days" column (M) takes its place visually.
$ws.Columns.Item(12).Hidden = $true
$ws.Columns.Item(13).ColumnWidth = 22.5

# --- Header for the new column ---
$ws.Range("M3").Value = "Revolution in Earth Days"

# --- Give the new column the same number format as column L (#,##0.00) ---
$ws.Range("M4:M14").NumberFormat = $ws.Range("L4").NumberFormat()

# Row 4 (Sun) - no revolution data, leave the new cell blank (style only)
# (Nothing to write - the NumberFormat assignment above already gave it style "4")

# Row 5 (Mercury) - single formula, not part of the fill-down shared group
$ws.Range("M5").Formula = "=L5*365"

# Rows 6-13 (Venus .. Uranus) - filled down as a shared formula
$ws.Range("M6:M13").Formula = "=L6*365"

# Row 14 (Moon) - literal value (27.32 days), not a formula
$ws.Range("M14").Value = 27.32

# --- New rows 19-20: scale-model reference values ---
$ws.Range("C19").Value = 6380
$ws.Range("C19").NumberFormat = $ws.Range("C4").NumberFormat()
$ws.Range("E19").Value = 1740
$ws.Range("E19").NumberFormat = $ws.Range("C4").NumberFormat()
$ws.Range("F19").Value = 384400

$ws.Range("C20").Value = 1.8
$ws.Range("E20").Value = 0.5
$ws.Range("F20").Value = 110.459770115

# --- Update the view state to match the author's working position ---
$ws.Application.ActiveWindow.Zoom = 136
$ws.Range("M18").Select()
